$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = -7.887
$ws.Range("C12").Value = -12.977
$ws.Range("D12").Value = -7.972
$ws.Range("D14").Value = -8.263
$ws.Range("D22").Value = -8.106
